$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Colossal sharp/stab"
$ws.Range("F1").Value = "Colossal blunt"
$ws.Range("E2").Value = "Colossal swords,`nColossal lances/spears"
$ws.Range("F2").Value = "Colossal hammers"

$ws.Range("D3:E3").Style = "Good"
$ws.Range("D4:E4").Style = "Good"

$ws.Range("I3").Select()
